$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the results matrix (row 4: RLBP, row 5: (19), row 6: (20), row 7: (13)) ---
# Row 4
$ws.Range("K4").Value = 0.035714285714285698
$ws.Range("L4").Value = 0.035714285714285698
$ws.Range("M4").Value = 0.32718750000000002
$ws.Range("N4").Value = 0.678448660714285

# Row 5
$ws.Range("K5").Value = 0.035714285714285698
$ws.Range("L5").Value = 0.035714285714285698
$ws.Range("M5").Value = 0.435926339285714
$ws.Range("N5").Value = 0.52712053571428497

# Row 6
$ws.Range("K6").Value = 0.062053571428571402
$ws.Range("L6").Value = 0.055145089285714197
$ws.Range("M6").Value = 0.37840401785714201
$ws.Range("N6").Value = 0.18966517857142801

# Row 7 - these cells were blank (K7, M7, N7 had the unused "fillId3/border1" style,
# L7 already carried the "fillId4/border1" style used by the rest of the row).
# Bring K7/M7/N7 into line with L7's formatting before filling in the new values.
$ws.Range("L7").Copy()
$ws.Range("K7").PasteSpecial(-4122)
$ws.Range("L7").Copy()
$ws.Range("M7").PasteSpecial(-4122)
$ws.Range("L7").Copy()
$ws.Range("N7").PasteSpecial(-4122)

$ws.Range("K7").Value = 0.91092633928571398
$ws.Range("L7").Value = 0.78564732142857097
$ws.Range("M7").Value = 0.93553571428571403
$ws.Range("N7").Value = 0.93553571428571403

# --- Move the active selection from K7 to N7 ---
$ws.Range("N7").Select()
